$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2-7), for columns D, L, M, N, O, P, Q, R, S, T
# This reflects a reshuffle of the weekly price records across dates.

# Row 2
$ws.Range("D2").Value = 44330
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = "`$/caja 18 kilos granel"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 861
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44698
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16500
$ws.Range("Q3").Value = "`$/caja 18 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 917
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44316
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 17500
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17750
$ws.Range("Q4").Value = "`$/caja 16 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1109
$ws.Range("T4").Value = 16

# Row 5
$ws.Range("D5").Value = 44316
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("Q5").Value = "`$/caja 16 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 16

# Row 6
$ws.Range("D6").Value = 44334
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 12500
$ws.Range("Q6").Value = "`$/caja 12 kilos empedrada"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1042
$ws.Range("T6").Value = 12

# Row 7
$ws.Range("D7").Value = 44344
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 13500
$ws.Range("Q7").Value = "`$/caja 18 kilos granel"
$ws.Range("R7").Value = "Provincia de Curicó"
$ws.Range("S7").Value = 750
$ws.Range("T7").Value = 18
